# update buglist and add pesantren list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Sheet1!B26 : update bug text
#    "fitur report berdasarkan pondok pesantren" -> "fitur daftar pondok pesantren"
# ---------------------------------------------------------------------
$ws.Range("B26").Value = "fitur daftar pondok pesantren"

# ---------------------------------------------------------------------
# 2) Row 24 : status moves from "open" to "close" (C24), reporter stays
#    "agung" (D24) and a new solver "agung" is recorded (E24)
#    (copy the "close" formatting - yellow fill - from an existing
#    closed-status cell so the style matches exactly)
# ---------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "close"

$ws.Range("D24").Value = "agung"
$ws.Range("E24").Value = "agung"

# ---------------------------------------------------------------------
# 3) Row 26 : same pattern as row 24
# ---------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "close"

$ws.Range("D26").Value = "agung"
$ws.Range("E26").Value = "agung"

# ---------------------------------------------------------------------
# 4) Row 28 : blank spacer row -- drop the leftover empty/styled cells
#    in C:E so the row only keeps the (also empty) A/B cells
# ---------------------------------------------------------------------
$ws.Range("C28:E28").Clear()

# ---------------------------------------------------------------------
# 5) View state : active selection moves to F9, scrolled back to top,
#    and the sheet-tab-area/horizontal-scroll-bar divider (tab ratio)
#    is resized
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.TabRatio = 12
$ws.Range("F9").Select()
